$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (style matches the existing header row, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Column I is a constant 1 for every data row; column J mirrors column H.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value()
}
